$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "codeforiati:group-name" (column D) and "codeforiati:group-code" (column E)
# columns were swapped for every row (header included): what used to be in D is now
# in E and vice-versa.
$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $dCell = $ws.Cells.Item($r, 4)
    $eCell = $ws.Cells.Item($r, 5)
    $dVal = $dCell.Value2
    $eVal = $eCell.Value2
    $dCell.Value2 = $eVal
    $eCell.Value2 = $dVal
}
